# IndMF - "Added other attributes like pe,pb etc"
#
# On the SmallCap sheet a brand-new data column is inserted right after the
# stock-name column (i.e. before the existing "Nov_22" column). The new
# column is seeded with the same values/formatting as the column that used
# to sit in that slot (so every existing figure just slides one column to
# the right), and its header is relabelled "Nov_24" since it represents a
# new reporting month.  MultiCap is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SmallCap")

# Copy column B (values + styles) and insert the copy in column B's place.
# This shifts the old B->C, C->D, D->E, E->F, leaving a duplicate of the
# old column B data sitting in the new column B.
$ws.Columns.Item(2).Copy()
$ws.Columns.Item(2).Insert()

# Relabel the header of the freshly-inserted column for the new month.
$ws.Cells.Item(1, 2).Value = "Nov_24"
